$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.776.57"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.96%  "
$ws.Range("D3").Value = "'1.865.43"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  +3.30%  "
$ws.Range("D5").Value = "'324.80"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.82%  "
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("E7").Value = "  +3.21%  "
$ws.Range("D8").Value = "'0.3803"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.77%  "
$ws.Range("D9").Value = "'0.07472"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("D10").Value = "'0.8855"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.73%  "
$ws.Range("D11").Value = "'21.73"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "'1.885.21"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -13.84%  "
$ws.Range("D13").Value = "'5.566"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").Value = "'6.768"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").Value = "'0.07240"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").Value = "'83.81"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").Value = "'1.040"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.35%  "
$ws.Range("D18").Value = "'0.000009167"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("E19").Value = "  +3.04%  "
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").Value = "'27.769.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("D22").Value = "'5.324"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("D23").Value = "'11.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.36%  "
$ws.Range("D24").Value = "'1.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +6.36%  "
$ws.Range("D25").Value = "'159.13"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("D26").Value = "'18.89"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.17%  "
$ws.Range("D27").Value = "'5.337"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("D28").Value = "'1.983"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.49%  "
$ws.Range("D29").Value = "'117.88"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("D30").Value = "'0.09072"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").Value = "'0.7799"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.38%  "
$ws.Range("D32").Value = "'3.104"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +10.53%  "
$ws.Range("D33").Value = "'1.214"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").Value = "'4.577"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.69%  "
$ws.Range("E35").Value = "  +3.14%  "
$ws.Range("D37").Value = "'0.01995"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.80%  "
$ws.Range("D38").Value = "'0.05352"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.88%  "
$ws.Range("D39").Value = "'2.865"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.25%  "
$ws.Range("D40").Value = "'0.5207"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("D42").Value = "'6.904"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.52%  "
$ws.Range("D43").Value = "'8.684"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("D44").Value = "'110.28"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.39%  "
$ws.Range("D45").Value = "'10.75"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.19%  "
$ws.Range("D46").Value = "'1.725"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.08%  "
$ws.Range("D47").Value = "'0.4714"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("D48").Value = "'0.06478"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.34%  "
$ws.Range("D49").Value = "'1.921"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.40%  "
$ws.Range("D50").Value = "'39.90"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("D51").Value = "'64.66"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.82%  "
